$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (ParticipantsTab): replace the Participant query (B2) with the
#     new, expanded query that adds OPTIONAL MATCH clauses and apoc.coll.sort
$participantQuery = @'
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE g.instrument_model in ['Illumina HiSeq 2000']
WITH p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN 
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samp, ','), '') as `Samples`
ORDER BY p.participant_id LIMIT 100
'@

# --- Row 3 (SamplesTab): wrap the tumor-status collection with apoc.coll.sort
$sampleQuery = @'
Match (f)<--(g:genomic_info)
WHERE g.instrument_model in ['Illumina HiSeq 2000']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p,s,samp,apoc.coll.sort(COLLECT(DISTINCT samp.sample_tumor_status)) as tumor
RETURN  
 coalesce(samp.sample_id, '') as `Sample ID`,
 coalesce(p.participant_id,'') as `Participant ID`,
 coalesce(s.study_name, '') as `Study Name`,
 coalesce(s.phs_accession,'') as `Accession`,
coalesce(samp.sample_tumor_status,'') as `Tumor`,
coalesce(samp.sample_type,'') as `Analyte Type`
ORDER By samp.sample_id LIMIT 100
'@

# --- Row 4 (FilesTab): wrap the tumor-status collection with apoc.coll.sort
$fileQuery = @'
Match (f)<--(g:genomic_info)
WHERE g.instrument_model in ['Illumina HiSeq 2000']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p,s,f,samp,apoc.coll.sort(COLLECT(DISTINCT samp.sample_tumor_status)) as tumor
RETURN 
    coalesce(f.file_name, '') as `File Name`,
    coalesce(s.study_name, '') as `Study Name`,
    coalesce(s.phs_accession,'') as `Accession`,
    coalesce(p.participant_id,'') as `Participant ID`,
    coalesce(samp.sample_id, '') as `Sample ID`,
    coalesce(f.file_type, '') as `File Type`
ORDER By f.file_name LIMIT 100
'@

$ws.Range("B2").Value = $participantQuery
$ws.Range("B3").Value = $sampleQuery
$ws.Range("B4").Value = $fileQuery

# Row heights grew because the replacement text is longer (wraps to more
# lines at the sheet's font/column width). Row 4's text length is unchanged
# so its height stays as-is.
$ws.Rows.Item(2).RowHeight = 354.75
$ws.Rows.Item(3).RowHeight = 282.75

# Update the view: no more scrolled top-left cell, and the active selection
# moves to C2 instead of D3.
$ws.Range("C2").Select() | Out-Null
